# Insert a new data row at row 472 (pushing existing rows 472:544 down to
# 473:545) and populate the new row with the latest price-report entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 472..544 down by one row.
$ws.Rows("472:472").Insert()

# Populate the newly inserted row 472 with the new record.
$ws.Range("A472").Value = 10
$ws.Range("B472").Value = "Vega Modelo de Temuco"
$ws.Range("C472").Value = "La Araucanía"
$ws.Range("D472").Value = 45218
$ws.Range("E472").Value = 9
$ws.Range("F472").Value = 100114013
$ws.Range("G472").Value = "Zanahoria"
$ws.Range("H472").Value = "Sin especificar"
$ws.Range("I472").Value = "Primera"
$ws.Range("J472").Value = 380
$ws.Range("K472").Value = 7000
$ws.Range("L472").Value = 7000
$ws.Range("M472").Value = 7000
$ws.Range("N472").Value = "$/saco 20 kilos"
$ws.Range("O472").Value = "Región Metropolitana"
$ws.Range("P472").Value = 350
$ws.Range("Q472").Value = 20
$ws.Range("R472").Value = "Hortaliza"
